$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) PLOG0014 subprocess title: retitle from
#    "PLOG0014 - Delimitação de valor de referência de material nacionalizado"
#    to
#    "PLOG0014 - Valor de referência a partir de histórico de requisições"
#    The final text is split across four runs (matching the source edit),
#    even though all four share identical run formatting.
# ---------------------------------------------------------------------------

$old_tail = "Delimitação de valor de referência de material nacionalizado"
$new_tail = "Valor de referência a partir de histórico de requisições"

$findRange = $d.Content
$replaced = $findRange.Find.Execute($old_tail, $true, $false, $false, $false, $false, $true, 1, $false, $new_tail, 1)

if ($replaced) {
    $locateRange = $d.Content
    $located = $locateRange.Find.Execute("PLOG0014 – " + $new_tail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($located) {
        $base = $locateRange.Start

        # Boundaries (relative to $base) splitting the sentence into:
        #   [0,11)  "PLOG0014 – "
        #   [11,12) "V"
        #   [12,31) "alor de referência "
        #   [31,67) "a partir de histórico de requisições"
        $segV = $d.Range($base + 11, $base + 12)
        $segV.Font.Bold = 1
        $segV.Font.Bold = 0

        $segAlor = $d.Range($base + 12, $base + 31)
        $segAlor.Font.Bold = 1
        $segAlor.Font.Bold = 0

        $segRest = $d.Range($base + 31, $base + 67)
        $segRest.Font.Bold = 1
        $segRest.Font.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# 2) Footer page-number field: the cached PAGE field result on the footer
#    used by the document's final section ("Página 4 de 7") is refreshed to
#    "Página 7 de 7".
# ---------------------------------------------------------------------------

$lastSection = $d.Sections.Item($d.Sections.Count)
$pageFooter = $lastSection.Footers.Item(1)
$pageUpdated = $pageFooter.Range.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "7", 1)
